# Fixes for watch list test cases
# Adds two new test-case rows to the "Test Cases" sheet:
#   - CommentsTabTimeStampValidationTest
#   - HCRProfileBadgeTest

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Copy the formatting of the last existing data row (row 37) down into the
# two new rows (38:39) so the new rows get the same cell styles (borders /
# fills / fonts) as the rest of the table.
$ws.Range("A37:E37").Copy() | Out-Null
$ws.Range("A38:E39").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Application.CutCopyMode = $false

# Row 38 - CommentsTabTimeStampValidationTest
$ws.Range("A38").Value = "CommentsTabTimeStampValidationTest"
$ws.Range("B38").Value = "TBD"
$ws.Range("C38").Value = "Verify that Comments tab comments displayed with timestamp"
$ws.Range("D38").Value = "Y"
$ws.Range("E38").Value = "SKIP"

# Row 39 - HCRProfileBadgeTest
$ws.Range("A39").Value = "HCRProfileBadgeTest"
$ws.Range("B39").Value = "TBD"
$ws.Range("C39").Value = "Verify that HCR profile having badge along with their name"
$ws.Range("D39").Value = "Y"
$ws.Range("E39").Value = "PASS"

# Update the sheet view to match the scrolled / selected state recorded in
# the workbook (top visible row 13, active cell C34).
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 13
$win.ScrollColumn = 1
$ws.Range("C34").Select() | Out-Null
